# Fruta / hortaliza, semanal
# Insert two new weekly price rows (393-394) for "Feria Lagunitas de Puerto
# Montt - Limon", pushing the existing rows 393..415 down to 395..417.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the position of the current row 393, shifting the
# remaining data (old rows 393-415) down to rows 395-417.
$ws.Range("A393:A394").EntireRow.Insert()

# Row 393 - "1a plateado", Provincia de Melipilla, $/malla 18 kilos
$ws.Cells.Item(393, 1).Value  = 4
$ws.Cells.Item(393, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(393, 3).Value  = "Los Lagos"
$ws.Cells.Item(393, 4).Value  = 44610
$ws.Cells.Item(393, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(393, 5).Value  = 10
$ws.Cells.Item(393, 6).Value  = "Fruta"
$ws.Cells.Item(393, 7).Value  = 100102
$ws.Cells.Item(393, 8).Value  = "Cítricos"
$ws.Cells.Item(393, 9).Value  = 100102003
$ws.Cells.Item(393, 10).Value = "Limón"
$ws.Cells.Item(393, 11).Value = "Sin especificar"
$ws.Cells.Item(393, 12).Value = "1a plateado"
$ws.Cells.Item(393, 13).Value = 1000
$ws.Cells.Item(393, 14).Value = 25000
$ws.Cells.Item(393, 15).Value = 25000
$ws.Cells.Item(393, 16).Value = 25000
$ws.Cells.Item(393, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(393, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(393, 19).Value = 1389
$ws.Cells.Item(393, 20).Value = 18

# Row 394 - "2a plateado", Provincia de Melipilla, $/malla 18 kilos
$ws.Cells.Item(394, 1).Value  = 4
$ws.Cells.Item(394, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(394, 3).Value  = "Los Lagos"
$ws.Cells.Item(394, 4).Value  = 44610
$ws.Cells.Item(394, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(394, 5).Value  = 10
$ws.Cells.Item(394, 6).Value  = "Fruta"
$ws.Cells.Item(394, 7).Value  = 100102
$ws.Cells.Item(394, 8).Value  = "Cítricos"
$ws.Cells.Item(394, 9).Value  = 100102003
$ws.Cells.Item(394, 10).Value = "Limón"
$ws.Cells.Item(394, 11).Value = "Sin especificar"
$ws.Cells.Item(394, 12).Value = "2a plateado"
$ws.Cells.Item(394, 13).Value = 500
$ws.Cells.Item(394, 14).Value = 21000
$ws.Cells.Item(394, 15).Value = 21000
$ws.Cells.Item(394, 16).Value = 21000
$ws.Cells.Item(394, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(394, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(394, 19).Value = 1167
$ws.Cells.Item(394, 20).Value = 18
